$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("E","F","G","H","I","J","K","L")

# Original E:L numeric values for rows 2..109, keyed by row number.
$data = @{}
$data[2] = @(16901.3286576178, 40942.4701175864, 115008.04067282, 104658.277264856, 250974.206737793, 648712.650075603, 190587.844116008, 0.603438488987871)
$data[3] = @(12930.8539166195, 20360.1418391378, 41327.0179412038, 38609.8247436038, 78329.2250881021, 142740.022550705, 190587.844116008, 0.216839736725542)
$data[4] = @(9277.57419635958, 17846.5193707543, 34252.7855019843, 32619.7433299421, 59842.7100761047, 133407.555966138, 190587.844116008, 0.179721774286587)
$data[5] = @(1517.80992042517, 4034.45480435485, 11350.4484034091, 10332.1577364619, 24890.1820481969, 72676.569904901, 19091.7805259841, 0.594520159498013)
$data[6] = @(950.549415310775, 1742.4676748416, 3577.22961736678, 3335.46814772499, 6746.69198814801, 13762.790464702, 19091.7805259841, 0.18737014143328)
$data[7] = @(1136.62571402095, 2094.6850606672, 4164.10250520819, 3944.3054155585, 7519.69966786055, 14250.9892109111, 19091.7805259841, 0.218109699068707)
$data[8] = @(5.41224095830597, 20.1714544381016, 79.6092497950228, 70.1768001444611, 194.184988908031, 461.750668690582, 134.744724504121, 0.590815336837831)
$data[9] = @(2.92604611993333, 6.33012491731484, 24.9622903334982, 18.8780863310735, 77.7177834972177, 196.304273970165, 134.744724504121, 0.185256160679855)
$data[10] = @(5.90056509257507, 12.0210339466061, 30.1731843755997, 27.7943508232212, 61.6113024335281, 118.644690862391, 134.744724504121, 0.223928502482314)
$data[11] = @(1.51591322079933, 5.78244178318538, 24.3093647962345, 21.3235267304314, 59.6238556030713, 185.0986111163, 49.9493225001049, 0.486680571016423)
$data[12] = @(4.18696810548222, 7.27809118971385, 19.9797297562294, 17.7613936084334, 45.7191473676192, 105.48789123351, 49.9493225001049, 0.400000015139093)
$data[13] = @(1.04513980887023, 2.2325183148282, 5.66022794764102, 5.13490379931018, 12.0739254577302, 32.5630695450657, 49.9493225001049, 0.113319413844485)
$data[14] = @(0.129365257350591, 0.606627803186471, 2.2287574177104, 1.95979577923768, 5.45254826027295, 13.0773748649964, 4.13906727296547, 0.538468517355986)
$data[15] = @(0.153414277243866, 0.303952100433841, 0.696816783871199, 0.62644221510576, 1.5092905597552, 3.16908871268211, 4.13906727296547, 0.168351161727303)
$data[16] = @(0.191578987323574, 0.393104363667808, 1.21349307138387, 1.0568310026981, 2.8603804268228, 7.16591118208625, 4.13906727296547, 0.29318032091671)
$data[17] = @(0.243810169778453, 0.839164164342427, 3.31002432141764, 2.90815982061156, 8.16480378789897, 20.4372323668144, 5.39021730985312, 0.614079939850854)
$data[18] = @(0.204717221697868, 0.493087278836272, 1.11807877816842, 1.02430797726613, 2.30465049725298, 4.65981326164781, 5.39021730985312, 0.207427402996278)
$data[19] = @(0.21584874754311, 0.391161492042017, 0.962114210267065, 0.853852244060294, 2.20155924136608, 5.1376096180102, 5.39021730985312, 0.178492657152868)
$data[20] = @(5.51926343296149, 15.7181314017232, 62.6289701530496, 55.0683301505356, 153.419233878511, 333.308906741706, 109.976969570647, 0.56947350338489)
$data[21] = @(8.13649171848497, 12.7594420073476, 26.5248940304697, 24.7978795622396, 50.5338283542505, 107.82158458723, 109.976969570647, 0.241185896774785)
$data[22] = @(5.07317657184204, 9.18937658725735, 20.8231053871274, 19.3317542872053, 40.9741592096047, 76.5533785453819, 109.976969570647, 0.189340599840325)
$data[23] = @(0.00872102170110465, 0.0259153377094936, 0.100807495734583, 0.0888534194107507, 0.244640006056124, 0.598278791963229, 0.184502187379284, 0.54637561302919)
$data[24] = @(0.00793185967321964, 0.0152888322387972, 0.0487716240542733, 0.0403292851533925, 0.132269778577808, 0.314161151272524, 0.184502187379284, 0.264341711862813)
$data[25] = @(0.00621650090193741, 0.012762589973438, 0.0349230675904277, 0.0314745732739767, 0.0761397945862771, 0.201450353653843, 0.184502187379284, 0.189282675107997)
$data[26] = @(1.86923918055446, 6.60724705478493, 28.5687758630722, 25.1798080825501, 70.6939907645111, 206.581127503203, 42.9372441346065, 0.665361190241047)
$data[27] = @(1.42294238559352, 3.61798169025557, 9.04391217151543, 8.18089229529873, 19.8105237293263, 46.107861628443, 42.9372441346065, 0.210630941826707)
$data[28] = @(1.35187562288018, 2.40647397962615, 5.3245561000189, 4.93058947021592, 10.5686435813594, 31.3174713177741, 42.9372441346065, 0.124007867932246)
$data[29] = @(1833.9972246051, 10939.5112904046, 25523.8415375821, 23845.2537324689, 49265.2542971955, 99702.4355794138, 72602.8387329154, 0.351554319128992)
$data[30] = @(1378.02718344846, 2339.78357337251, 4322.89389149724, 4097.5451673435, 7667.8340546165, 14973.7520689903, 72602.8387329154, 0.0595416648569335)
$data[31] = @(19201.0524226326, 27704.7268490638, 42756.1033038361, 41917.6954688038, 62899.8669787849, 89288.9699524989, 72602.8387329154, 0.588904016014074)
$data[32] = @(242.937118100356, 1334.93174978624, 3229.29562205382, 2991.39628558266, 6501.73286122857, 12606.882692223, 8989.67980945146, 0.359222540791568)
$data[33] = @(109.055672410561, 202.936336779937, 376.136188413106, 356.525831624318, 663.20390945203, 1471.75750383306, 8989.67980945146, 0.0418408882614093)
$data[34] = @(2169.90104464188, 3257.09338394837, 5384.24799898454, 5229.37856477796, 8372.02360290906, 12325.9272087131, 8989.67980945146, 0.598936570947022)
$data[35] = @(1.42573460464363, 6.33803656908018, 20.8309876992354, 18.7527451342416, 47.5024018682133, 104.048916658654, 54.7652070848, 0.380369011788454)
$data[36] = @(0.473871958443119, 0.906496661736256, 2.74373528825468, 2.28201029090451, 7.10304887256855, 16.0233342359836, 54.7652070848, 0.0500999710273387)
$data[37] = @(9.14022181309656, 15.7535133109257, 31.1904840973099, 29.7982480952893, 54.4529924510039, 96.5266005115867, 54.7652070848, 0.569531017184207)
$data[38] = @(0.0915138717998519, 0.388302869040905, 1.21716644683402, 1.06105555851873, 2.97801298002335, 12.643902806589, 5.7666566058659, 0.211069694282803)
$data[39] = @(0.532568373492503, 1.03171161268866, 2.33696856754407, 2.1417493907169, 4.68558068523006, 11.203621055053, 5.7666566058659, 0.405255371919819)
$data[40] = @(0.730647530113912, 1.20295772450213, 2.21252159148781, 2.09400831431796, 3.91178290576803, 7.45171981055594, 5.7666566058659, 0.383674933797378)
$data[41] = @(0.0272116381929811, 0.137458242828243, 0.887978376594552, 0.680976748733618, 2.68047429766662, 6.3330583885319, 2.43316681295445, 0.364947594988908)
$data[42] = @(0.0209785667459499, 0.0398162749642096, 0.0860115038000885, 0.0794925680418175, 0.168526991016872, 0.388225266882341, 2.43316681295445, 0.0353496124236751)
$data[43] = @(0.242876794546294, 0.556941042767339, 1.45917693255981, 1.31284082203174, 3.07751987784068, 5.45447979988605, 2.43316681295445, 0.599702792587417)
$data[44] = @(0.0592516279268513, 0.20890439352135, 0.770910449336563, 0.638153790864921, 2.12608634253561, 6.43285001627765, 2.05165970771063, 0.375749665716638)
$data[45] = @(0.0323767493785898, 0.0595181953023992, 0.123006209993833, 0.115308843965431, 0.233132804034115, 0.479801852495671, 2.05165970771063, 0.0599544893003195)
$data[46] = @(0.364570440206111, 0.582881504080713, 1.15774304838023, 1.05920419688264, 2.31734358924616, 4.05920564137116, 2.05165970771063, 0.564295844983042)
$data[47] = @(0.999655960426012, 4.69008913602796, 15.5928265753228, 13.7815178024498, 36.4096973300956, 69.2589257054728, 40.9266705503669, 0.380994260359717)
$data[48] = @(0.871251060902835, 1.57608363332891, 2.99252595383701, 2.83908621805059, 5.3894569651022, 9.22056031592416, 40.9266705503669, 0.0731192133050312)
$data[49] = @(7.91158353775744, 11.7478369047042, 22.341318021207, 21.2568297833758, 39.3220426665407, 77.0729740562191, 40.9266705503669, 0.545886526335251)
$data[50] = @(0.00143884486117885, 0.00551851024996745, 0.0224532929906544, 0.0189633417465717, 0.060840542128592, 0.172409794023171, 0.0685609524692209, 0.327493889480815)
$data[51] = @(0.0013936932904714, 0.00281766918870271, 0.00812671918679765, 0.00719973671378439, 0.0185514678302715, 0.0447186370625334, 0.0685609524692209, 0.118532763827136)
$data[52] = @(0.0103860994481277, 0.018280241804221, 0.0379809402917688, 0.0353493036760269, 0.0734405140307679, 0.125283697563146, 0.0685609524692209, 0.553973346692049)
$data[53] = @(0.144652685517423, 0.908753131160276, 2.80791154514978, 2.42645957353023, 7.41075547988746, 21.3385402697675, 12.6934794065257, 0.221208973144609)
$data[54] = @(0.229438226303559, 0.464103773627472, 1.14222289081269, 1.02481497404547, 2.48565842899028, 4.57137303562461, 12.6934794065257, 0.0899850115347792)
$data[55] = @(1.82244564886408, 3.56506700636163, 8.74334497056319, 7.71173506359619, 18.7453158285, 32.1678415398826, 12.6934794065257, 0.688806015320612)
$data[56] = @(37485.8850021217, 82521.8325764481, 185742.415412402, 172220.247001479, 370535.116298305, 1162080.78676884, 315157.391143707, 0.589363983304793)
$data[57] = @(22585.4438934668, 52658.2341291492, 116427.92877195, 108612.704061631, 226796.145778592, 486790.83614605, 315157.391143707, 0.369427886014137)
$data[58] = @(2356.8141208799, 4695.55129349055, 12987.0469593547, 11653.0741574609, 29171.2957385752, 63541.9634823214, 315157.391143707, 0.0412081306810693)
$data[59] = @(4497.78767922878, 8571.49207625289, 19121.4791014892, 17809.4247063026, 37342.7963180421, 95912.8846540229, 30988.0410033223, 0.617059952238968)
$data[60] = @(1978.91023693746, 4521.17293686401, 10136.7650762886, 9437.69648415862, 19662.6705095946, 40309.2977873247, 30988.0410033223, 0.327118615700869)
$data[61] = @(250.254987309596, 568.523618451195, 1729.7968255445, 1497.23091928841, 4267.08195729896, 12431.2545534793, 30988.0410033223, 0.055821432060163)
$data[62] = @(17.7810446927995, 45.4823334279642, 132.140007100245, 120.079805374567, 288.984109282236, 668.009249130249, 212.163381607669, 0.622821931376441)
$data[63] = @(5.94351515482316, 18.3731605826389, 71.4255517783051, 56.5080789226577, 209.468631393121, 538.623157325886, 212.163381607669, 0.336653531995379)
$data[64] = @(0.930778448813074, 2.42313655028105, 8.59782272911849, 7.12793429702034, 23.8309439610768, 90.6976707327598, 212.163381607669, 0.0405245366281799)
$data[65] = @(3.09890713403993, 9.28871804547405, 34.7776240905542, 30.8313577384305, 82.3636772500776, 234.552731104892, 92.9663241604445, 0.374088406792698)
$data[66] = @(7.94470768553644, 20.4675211222981, 57.6371680277117, 51.5320775021249, 127.364562648516, 407.138571798427, 92.9663241604445, 0.619978992911879)
$data[67] = @(0.0925631779588281, 0.185411512275483, 0.551532042178589, 0.48090405232419, 1.31108088380479, 4.19212880214866, 92.9663241604445, 0.00593260029542242)
$data[68] = @(0.587192250515534, 1.27510728650804, 4.03891756301732, 3.67005880021712, 8.9365188058429, 21.1765185606762, 6.56013552016912, 0.615675933919304)
$data[69] = @(0.316483183112249, 0.805596182217462, 2.08393126075638, 1.88906571869131, 4.51491477450707, 12.5617687747501, 6.56013552016912, 0.31766588576552)
$data[70] = @(0.0471164925402921, 0.111797065614206, 0.437286696395419, 0.360118894838093, 1.2143093289787, 3.74009509032166, 6.56013552016912, 0.0666581803151753)
$data[71] = @(0.737262372484103, 1.73661515238413, 5.38165480249472, 4.82968720702311, 11.9666078802037, 45.2102868000248, 8.90516343303795, 0.604329706350913)
$data[72] = @(0.645008568351789, 1.30476750005855, 3.20259704783131, 2.93487254545793, 6.60611993872416, 14.5475722746166, 8.90516343303795, 0.359633719460975)
$data[73] = @(0.0407476766151235, 0.108769555690368, 0.32091158271193, 0.283168742430714, 0.755565068404437, 1.82858236577505, 8.90516343303795, 0.0360365741881115)
$data[74] = @(11.5839863346293, 34.2892271992753, 103.011457701133, 93.660635417557, 225.373981361668, 619.131449636669, 185.39419212405, 0.555634761374871)
$data[75] = @(15.6683744272472, 33.8031641928768, 76.2745769352973, 70.96126486838, 150.180595568295, 334.360017645357, 185.39419212405, 0.411418373258754)
$data[76] = @(1.00310776825539, 2.10648012125371, 6.10815748761901, 5.40089582681285, 14.4403195196493, 37.4688357026068, 185.39419212405, 0.0329468653663755)
$data[77] = @(0.0192498206205066, 0.0512285719135253, 0.163397834574768, 0.147982879476975, 0.362700602493682, 1.12913403039066, 0.33608643032231, 0.48617801801182)
$data[78] = @(0.0135896024530579, 0.0467681446846812, 0.16110290894438, 0.138419011690501, 0.406452248116863, 0.87556225062238, 0.33608643032231, 0.479349638692286)
$data[79] = @(0.00118282335236916, 0.00361844997320529, 0.0115856868031621, 0.00995707292877167, 0.0298740302942119, 0.0891957902068791, 0.33608643032231, 0.0344723432958936)
$data[80] = @(4.26649798503975, 12.0187262779557, 42.2462620135755, 37.5425490115745, 100.460633477657, 342.03476624524, 72.7981226934644, 0.580320761724372)
$data[81] = @(4.39628712226956, 9.97581268917254, 27.0013894963876, 24.1776406342893, 60.7933225377175, 188.588533483386, 72.7981226934644, 0.370907772032584)
$data[82] = @(0.348371565280823, 0.850562756086513, 3.55047118350124, 2.85285567447922, 10.074063279967, 24.8473687524947, 72.7981226934644, 0.0487714662430435)
$data[83] = @(20226.6728015455, 36445.2608986991, 72027.8110547414, 68110.9870508976, 130953.096143361, 230011.758731049, 109632.718183634, 0.656991929490386)
$data[84] = @(6221.24491181103, 12627.1074194985, 31300.8883505003, 28708.3836964068, 64217.8021768321, 128377.918054312, 109632.718183634, 0.285506816478559)
$data[85] = @(282.693753777517, 1348.23310245164, 6304.01877839222, 5021.47640263776, 18463.8270591297, 55265.9171451548, 109632.718183634, 0.0575012540310552)
$data[86] = @(1900.00818884544, 4030.49955636955, 8223.58252576977, 7710.60032672319, 15275.891465178, 26113.1538812114, 11627.4320163108, 0.707256986257488)
$data[87] = @(493.544582585718, 1091.60387288115, 2743.63337549617, 2518.77372320388, 5667.61722122611, 11405.4930546028, 11627.4320163108, 0.235962108541891)
$data[88] = @(26.0061043954443, 128.375467772468, 660.21611504481, 511.736659271085, 2033.77844345819, 7096.85897823526, 11627.4320163108, 0.0567809052006213)
$data[89] = @(8.35552200050403, 21.3965290776123, 55.0473520308457, 51.038673508104, 112.958795736667, 237.335442189214, 77.2492739241422, 0.712593779003043)
$data[90] = @(1.73127146395646, 4.76834903935571, 19.5380677542503, 15.1212895389697, 60.6481476889949, 168.282485683525, 77.2492739241422, 0.252922348156132)
$data[91] = @(0.153931032524107, 0.522752852957649, 2.66385413904625, 2.08678319082367, 8.18539714547642, 25.5160016593288, 77.2492739241422, 0.0344838728408259)
$data[92] = @(0.936281326890123, 2.82184615709978, 8.78593699666615, 7.98834029722014, 19.2068469575164, 75.4018914630803, 24.8630938167296, 0.35337263582034)
$data[93] = @(2.1097513429533, 5.25920472132413, 15.7828743631769, 13.7947738694034, 37.5383669052103, 79.3990026687811, 24.8630938167296, 0.634791248406788)
$data[94] = @(0.0133492983087749, 0.0568855837671703, 0.294282456886585, 0.225401762666181, 0.917752987424922, 3.91698468446754, 24.8630938167296, 0.0118361157728718)
$data[95] = @(0.181016821229249, 0.547587042419786, 2.01383877988226, 1.68759038006516, 5.28608673609628, 14.4929140304583, 2.76176169476177, 0.72918629572635)
$data[96] = @(0.102466483556536, 0.209410870311323, 0.59454882536393, 0.523901607972199, 1.39894702561226, 3.00606458458655, 2.76176169476177, 0.215278829629512)
$data[97] = @(0.00738244266560109, 0.0219335339947498, 0.153374089515578, 0.109699819295247, 0.54016027154362, 2.13874174096797, 2.76176169476177, 0.0555348746441383)
$data[98] = @(0.315520904455785, 0.781894467581296, 2.13354912395895, 1.91790634459177, 4.66180407997719, 12.5598404320148, 3.14931007223196, 0.677465563893135)
$data[99] = @(0.158118088844387, 0.325259208182205, 0.875318909316628, 0.785597460469518, 1.92804516162622, 4.47622029274471, 3.14931007223196, 0.277939894529432)
$data[100] = @(0.00599470877827365, 0.0275224889614463, 0.140442038956378, 0.110430690833284, 0.436009564453251, 1.60468426394375, 3.14931007223196, 0.0445945415774333)
$data[101] = @(6.39141539829279, 16.7484644810703, 42.0179715830034, 38.853516547982, 86.9019746981263, 181.963199070562, 65.2840660646109, 0.643617564221852)
$data[102] = @(4.14055860337006, 8.34472827447249, 20.889176656606, 19.136421137277, 43.1638067601502, 89.5661452987742, 65.2840660646109, 0.319973584916298)
$data[103] = @(0.130526545585815, 0.496427813789741, 2.37691782500157, 1.88185932073075, 7.1538332644113, 20.6270111025026, 65.2840660646109, 0.0364088508618498)
$data[104] = @(0.00978736757968176, 0.0228833743736855, 0.0634580933355305, 0.0574858379771684, 0.138600958598226, 0.308307343790889, 0.117306557729303, 0.540959470330438)
$data[105] = @(0.00538196461909998, 0.0134381461643416, 0.04948122127324, 0.0416097087520696, 0.132718737699118, 0.297537437370015, 0.117306557729303, 0.42181121184566)
$data[106] = @(0.000181337771700954, 0.00084055407634738, 0.00436724312053221, 0.00332893394544846, 0.0135532903214902, 0.055357928103214, 0.117306557729303, 0.0372293178239028)
$data[107] = @(1.54857325441624, 4.30622804794742, 12.329468528751, 11.3433111315058, 26.194290018302, 59.9167101478629, 22.0296809285362, 0.55967531117439)
$data[108] = @(1.07488022874132, 2.47476505877617, 7.62580673564281, 6.68791114938538, 18.2635803186816, 43.5301711693472, 22.0296809285362, 0.346160562215166)
$data[109] = @(0.0661776959871535, 0.24384497744814, 2.07440566414234, 1.34844805863833, 8.29622697253453, 35.3567404017832, 22.0296809285362, 0.0941641266104431)

# Rows 2-55 held "Northeast Atlantic" data; rows 56-109 held "Central North Atlantic" data.
# The edit swaps the data (and area label) between row i and row (i + 54).
for ($i = 2; $i -le 55; $i++) {
    $j = $i + 54
    $rowI = $data[$i]
    $rowJ = $data[$j]
    for ($k = 0; $k -lt $cols.Length; $k++) {
        $ws.Range($cols[$k] + $i).Value = $rowJ[$k]
        $ws.Range($cols[$k] + $j).Value = $rowI[$k]
    }
    $ws.Cells.Item($i, 1).Value = "Central North Atlantic"
    $ws.Cells.Item($j, 1).Value = "Northeast Atlantic"
}

